$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("F2").Value = "2021-10-05 13:39:34.117495"
$ws.Range("F3").Value = "2021-10-05 13:39:34.117507"
$ws.Range("F4").Value = "2021-10-05 13:39:34.117510"
$ws.Range("F5").Value = "2021-10-05 13:39:34.117513"
$ws.Range("F6").Value = "2021-10-05 13:39:34.117516"
$ws.Range("F7").Value = "2021-10-05 13:39:34.117519"
$ws.Range("F8").Value = "2021-10-05 13:39:34.117522"
$ws.Range("F9").Value = "2021-10-05 13:39:34.117524"
$ws.Range("F10").Value = "2021-10-05 13:39:34.117527"
$ws.Range("F11").Value = "2021-10-05 13:39:34.117530"
$ws.Range("F12").Value = "2021-10-05 13:39:34.117533"
$ws.Range("F13").Value = "2021-10-05 13:39:34.117535"
$ws.Range("F14").Value = "2021-10-05 13:39:34.117538"
$ws.Range("F15").Value = "2021-10-05 13:39:34.117540"
$ws.Range("F16").Value = "2021-10-05 13:39:34.117543"
$ws.Range("F17").Value = "2021-10-05 13:39:34.117546"
$ws.Range("F18").Value = "2021-10-05 13:39:34.117549"
$ws.Range("F19").Value = "2021-10-05 13:39:34.117552"
$ws.Range("F20").Value = "2021-10-05 13:39:34.117555"
$ws.Range("F21").Value = "2021-10-05 13:39:34.117557"
$ws.Range("F22").Value = "2021-10-05 13:39:34.117560"
$ws.Range("F23").Value = "2021-10-05 13:39:34.117562"
$ws.Range("F24").Value = "2021-10-05 13:39:34.117565"
$ws.Range("F25").Value = "2021-10-05 13:39:34.117568"
$ws.Range("F26").Value = "2021-10-05 13:39:34.117571"
$ws.Range("F27").Value = "2021-10-05 13:39:34.117573"
$ws.Range("F28").Value = "2021-10-05 13:39:34.117576"
$ws.Range("F29").Value = "2021-10-05 13:39:34.117578"
$ws.Range("F30").Value = "2021-10-05 13:39:34.117581"
$ws.Range("F31").Value = "2021-10-05 13:39:34.117584"

Write-Output "done"
